$wb = $excel.ActiveWorkbook

$wsActivos = $wb.Worksheets.Item("activos")
$wsPasivos = $wb.Worksheets.Item("pasivos")
$wsPatrimonio = $wb.Worksheets.Item("patrimonio")
$wsOtros = $wb.Worksheets.Item("otros")

# --- activos ---
$wsActivos.Range("A2").Value = "cuentas comerciales por cobrar y otras cuentas por cobrar no corrientes"
$wsActivos.Range("B2").NumberFormat = "@"
$wsActivos.Range("B2").Value = "138313"
$wsActivos.Range("B2").Style = "Normal"
$wsActivos.Range("A3").Value = "inversiones contabilizadas utilizando el método de la participación"
$wsActivos.Range("B3").NumberFormat = "@"
$wsActivos.Range("B3").Value = "331502901"
$wsActivos.Range("B3").Style = "Normal"
$wsActivos.Range("A4").Value = "Inversiones contabilizadas utilizando el método de la participación"
$wsActivos.Range("B4").NumberFormat = "@"
$wsActivos.Range("B4").Value = "334657003"
$wsActivos.Range("B4").Style = "Normal"
$wsActivos.Range("B7").NumberFormat = "@"
$wsActivos.Range("B7").Value = "4110735696"
$wsActivos.Range("B7").Style = "Normal"
$wsActivos.Range("B8").NumberFormat = "@"
$wsActivos.Range("B8").Value = "3392498650"
$wsActivos.Range("B8").Style = "Normal"
$wsActivos.Range("B9").NumberFormat = "@"
$wsActivos.Range("B9").Value = "369637458"
$wsActivos.Range("B9").Style = "Normal"

# --- pasivos ---
$wsPasivos.Range("B2").NumberFormat = "@"
$wsPasivos.Range("B2").Value = "2653580482"
$wsPasivos.Range("B2").Style = "Normal"
$wsPasivos.Range("B4").NumberFormat = "@"
$wsPasivos.Range("B4").Value = "3401565"
$wsPasivos.Range("B4").Style = "Normal"
$wsPasivos.Range("B5").NumberFormat = "@"
$wsPasivos.Range("B5").Value = "3263065"
$wsPasivos.Range("B5").Style = "Normal"

# --- patrimonio ---
$wsPatrimonio.Range("A5").Value = "otras reservas"
$wsPatrimonio.Range("B6").NumberFormat = "@"
$wsPatrimonio.Range("B6").Value = "633715769"
$wsPatrimonio.Range("B6").Style = "Normal"
$wsPatrimonio.Range("A7").Value = "Participaciones no controladoras"
$wsPatrimonio.Range("B7").NumberFormat = "@"
$wsPatrimonio.Range("B7").Value = "607015945"
$wsPatrimonio.Range("B7").Style = "Normal"

# --- otros ---
$wsOtros.Range("A11").Value = "Activos intangibles distintos de la plusvalía"
$wsOtros.Range("B11").NumberFormat = "@"
$wsOtros.Range("B11").Value = "862236570"
$wsOtros.Range("B11").Style = "Normal"
$wsOtros.Range("A12").Value = "Plusvalía"
$wsOtros.Range("B12").NumberFormat = "@"
$wsOtros.Range("B12").Value = "2059796230"
$wsOtros.Range("B12").Style = "Normal"
$wsOtros.Range("A13").Value = "planta y equipo"
$wsOtros.Range("B13").NumberFormat = "@"
$wsOtros.Range("B13").Value = "3743122719"
$wsOtros.Range("B13").Style = "Normal"
$wsOtros.Range("A14").Value = "Propiedad de inversión"
$wsOtros.Range("B14").NumberFormat = "@"
$wsOtros.Range("B14").Value = "3188927576"
$wsOtros.Range("B14").Style = "Normal"
$wsOtros.Range("A15").Value = "no corrientes"
$wsOtros.Range("B15").NumberFormat = "@"
$wsOtros.Range("B15").Value = "4046018"
$wsOtros.Range("B15").Style = "Normal"
$wsOtros.Range("A16").Value = "Activos por impuestos diferidos"
$wsOtros.Range("B16").NumberFormat = "@"
$wsOtros.Range("B16").Value = "356550480"
$wsOtros.Range("B16").Style = "Normal"
$wsOtros.Range("A17").Value = "MS PASIVOS CORRIENTES Otros pasivos financieros corrientes"
$wsOtros.Range("B17").NumberFormat = "@"
$wsOtros.Range("B17").Value = "505461062"
$wsOtros.Range("B17").Style = "Normal"
$wsOtros.Range("A18").Value = "Cuentas por pagar comerciales y otras cuentas por pagar"
$wsOtros.Range("B18").NumberFormat = "@"
$wsOtros.Range("B18").Value = "2866975457"
$wsOtros.Range("B18").Style = "Normal"
$wsOtros.Range("A19").Value = "Otras provisiones corrientes"
$wsOtros.Range("B19").NumberFormat = "@"
$wsOtros.Range("B19").Value = "16826672"
$wsOtros.Range("B19").Style = "Normal"
$wsOtros.Range("A20").Value = "Provisiones corrientes por beneficios a los empleados"
$wsOtros.Range("B20").NumberFormat = "@"
$wsOtros.Range("B20").Value = "130178251"
$wsOtros.Range("B20").Style = "Normal"
$wsOtros.Range("A21").Value = "Otros pasivos no financieros corrientes"
$wsOtros.Range("B21").NumberFormat = "@"
$wsOtros.Range("B21").Value = "240505744"
$wsOtros.Range("B21").Style = "Normal"
$wsOtros.Range("A22").Value = "PASIVOS NO CORRIENTES Otros pasivos financieros no corrientes"
$wsOtros.Range("B22").NumberFormat = "@"
$wsOtros.Range("B22").Value = "3704831700"
$wsOtros.Range("B22").Style = "Normal"
$wsOtros.Range("A23").Value = "Pasivos por arrendamientos no corrientes"
$wsOtros.Range("B23").NumberFormat = "@"
$wsOtros.Range("B23").Value = "1098575638"
$wsOtros.Range("B23").Style = "Normal"
$wsOtros.Range("A24").Value = "Cuentas comerciales por pagar y otras cuentas por pagar no corrientes"
$wsOtros.Range("B24").NumberFormat = "@"
$wsOtros.Range("B24").Value = "3536289"
$wsOtros.Range("B24").Style = "Normal"
$wsOtros.Range("A25").Value = "Otras provisiones no corrientes"
$wsOtros.Range("B25").NumberFormat = "@"
$wsOtros.Range("B25").Value = "48070186"
$wsOtros.Range("B25").Style = "Normal"
$wsOtros.Range("A26").Value = "Pasivo por impuestos diferidos"
$wsOtros.Range("B26").NumberFormat = "@"
$wsOtros.Range("B26").Value = "558350832"
$wsOtros.Range("B26").Style = "Normal"
$wsOtros.Range("A27").Value = "Otros pasivos no financieros no corrientes"
$wsOtros.Range("B27").NumberFormat = "@"
$wsOtros.Range("B27").Value = "76027357"
$wsOtros.Range("B27").Style = "Normal"
$wsOtros.Range("A28").Value = "MS MS Ingresos de actividades ordinarias"
$wsOtros.Range("B28").NumberFormat = "@"
$wsOtros.Range("B28").Value = "3503183757"
$wsOtros.Range("B28").Style = "Normal"
$wsOtros.Range("A29").Value = "Otros ingresos"
$wsOtros.Range("B29").NumberFormat = "@"
$wsOtros.Range("B29").Value = "8518912"
$wsOtros.Range("B29").Style = "Normal"
$wsOtros.Range("A30").Value = "Ingresos financieros"
$wsOtros.Range("B30").NumberFormat = "@"
$wsOtros.Range("B30").Value = "4602195"
$wsOtros.Range("B30").Style = "Normal"
$wsOtros.Range("A31").Value = "Efectivo y equivalentes al efectivo al principio del período"
$wsOtros.Range("B31").NumberFormat = "@"
$wsOtros.Range("B31").Value = "373700303"
$wsOtros.Range("B31").Style = "Normal"
$wsOtros.Range("A32").Value = "Efectivo y equivalentes al efectivo al final del período"
$wsOtros.Range("B32").NumberFormat = "@"
$wsOtros.Range("B32").Value = "564926038"
$wsOtros.Range("B32").Style = "Normal"
